$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-16 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-17 Sunday", 2)

$d.Content.Find.Execute("98×14=", $true, $false, $false, $false, $false, $true, 1, $false, "97×96=", 2)
$d.Content.Find.Execute("50×49=", $true, $false, $false, $false, $false, $true, 1, $false, "62×17=", 2)
$d.Content.Find.Execute("71×69=", $true, $false, $false, $false, $false, $true, 1, $false, "98×55=", 2)
$d.Content.Find.Execute("45×99=", $true, $false, $false, $false, $false, $true, 1, $false, "59×93=", 2)
$d.Content.Find.Execute("46×84=", $true, $false, $false, $false, $false, $true, 1, $false, "17×81=", 2)
$d.Content.Find.Execute("78×63=", $true, $false, $false, $false, $false, $true, 1, $false, "19×36=", 2)
$d.Content.Find.Execute("45×95=", $true, $false, $false, $false, $false, $true, 1, $false, "86×27=", 2)
$d.Content.Find.Execute("80×15=", $true, $false, $false, $false, $false, $true, 1, $false, "17×69=", 2)
$d.Content.Find.Execute("76×30=", $true, $false, $false, $false, $false, $true, 1, $false, "60×39=", 2)
$d.Content.Find.Execute("74×36=", $true, $false, $false, $false, $false, $true, 1, $false, "67×29=", 2)
$d.Content.Find.Execute("19×74=", $true, $false, $false, $false, $false, $true, 1, $false, "28×49=", 2)
$d.Content.Find.Execute("49×53=", $true, $false, $false, $false, $false, $true, 1, $false, "77×81=", 2)
$d.Content.Find.Execute("34×38=", $true, $false, $false, $false, $false, $true, 1, $false, "93×15=", 2)
$d.Content.Find.Execute("14×41=", $true, $false, $false, $false, $false, $true, 1, $false, "57×48=", 2)
$d.Content.Find.Execute("52×78=", $true, $false, $false, $false, $false, $true, 1, $false, "80×86=", 2)
$d.Content.Find.Execute("72×95=", $true, $false, $false, $false, $false, $true, 1, $false, "77×61=", 2)
$d.Content.Find.Execute("96×80=", $true, $false, $false, $false, $false, $true, 1, $false, "22×19=", 2)
$d.Content.Find.Execute("55×99=", $true, $false, $false, $false, $false, $true, 1, $false, "37×70=", 2)
$d.Content.Find.Execute("28×72=", $true, $false, $false, $false, $false, $true, 1, $false, "43×88=", 2)
$d.Content.Find.Execute("60×22=", $true, $false, $false, $false, $false, $true, 1, $false, "45×63=", 2)
$d.Content.Find.Execute("47×25=", $true, $false, $false, $false, $false, $true, 1, $false, "48×40=", 2)
$d.Content.Find.Execute("99×95=", $true, $false, $false, $false, $false, $true, 1, $false, "51×62=", 2)
$d.Content.Find.Execute("69×81=", $true, $false, $false, $false, $false, $true, 1, $false, "84×88=", 2)
$d.Content.Find.Execute("40×39=", $true, $false, $false, $false, $false, $true, 1, $false, "65×85=", 2)
$d.Content.Find.Execute("16×26=", $true, $false, $false, $false, $false, $true, 1, $false, "89×51=", 2)
